$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.344.37'
$ws.Range('E2').Value = '  -2.67%  '
$ws.Range('D3').Value = '3.626.59'
$ws.Range('E3').Value = '  +0.51%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.78'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.99'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -4.84%  '
$ws.Range('D7').Value = '3.617.14'
$ws.Range('E7').Value = '  +0.48%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.607'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.23%  '
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('E10').Value = '  -5.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.98'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +22.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.604'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.87%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '48.10'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -4.13%  '
$ws.Range('D15').Value = '4.209.16'
$ws.Range('E15').Value = '  +0.45%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '674.06'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -3.27%  '
$ws.Range('E17').Value = '  -1.00%  '
$ws.Range('D18').Value = '3.631.01'
$ws.Range('E18').Value = '  +0.69%  '
$ws.Range('D19').Value = '70.405.92'
$ws.Range('E19').Value = '  -2.65%  '
$ws.Range('E20').Value = '  -0.24%  '
$ws.Range('E21').Value = '  -4.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.40'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.936'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '17.01'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '99.62'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -5.22%  '
$ws.Range('E26').Value = '  -3.39%  '
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('E29').Value = '  -3.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.47'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.98'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.98%  '
$ws.Range('E32').Value = '  -5.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.49'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('E34').Value = '  -6.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.94'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -5.27%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '578.00'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -3.00%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '11.02'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.85%  '
$ws.Range('E38').Value = '  -0.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '58.14'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.28%  '
$ws.Range('E40').Value = '  +0.06%  '
$ws.Range('D41').Value = '3.551.25'
$ws.Range('E41').Value = '  -3.13%  '
$ws.Range('E43').Value = '  -3.20%  '
$ws.Range('E44').Value = '  -1.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '34.40'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -5.01%  '
$ws.Range('D46').Value = '0.0₃0728'
$ws.Range('E46').Value = '  -7.56%  '
$ws.Range('E47').Value = '  -6.05%  '
$ws.Range('E48').Value = '  +1.88%  '
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '136.35'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.90%  '
$ws.Range('E51').Value = '  -1.32%  '
